$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Name", "Course", "Semester", "Form Number", "Contact Number", "Email id", "Address")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Cells.Item(2, 6).NumberFormat = "@"
$row2 = @("Aaakef", "kjkjl", "kjkjk", "nnnnnn", "mmmm", "111", "mm333")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$ws.Columns.Item(1).ColumnWidth = 29.1666666666667
$ws.Columns.Item(2).ColumnWidth = 9.16666666666667
$ws.Columns.Item(3).ColumnWidth = 9.16666666666667
$ws.Columns.Item(4).ColumnWidth = 19.1666666666667
$ws.Columns.Item(5).ColumnWidth = 19.1666666666667
$ws.Columns.Item(6).ColumnWidth = 39.1666666666667
$ws.Columns.Item(7).ColumnWidth = 49.1666666666667
